$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.888.21"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "3.326.52"
$ws.Range("E3").Value = "  +0.86%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.43"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.58%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "3.323.47"
$ws.Range("E9").Value = "  +1.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.578"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.00"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.09%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.82%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "697.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.84%  "

# Row 15
$ws.Range("D15").Value = "3.870.42"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
$ws.Range("D17").Value = "67.899.16"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("E18").Value = "  -0.95%  "

# Row 19
$ws.Range("D19").Value = "3.335.44"
$ws.Range("E19").Value = "  +0.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.97"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.891"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.67%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.93%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.34%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.73"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.69%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.66%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.36"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.52"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.35%  "

# Row 31
$ws.Range("E31").Value = "  +7.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.89"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.97"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.26%  "

# Row 34
$ws.Range("E34").Value = "  +0.91%  "

# Row 36
$ws.Range("D36").Value = "3.708.80"
$ws.Range("E36").Value = "  -5.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.12"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.33%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.20%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.59"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.80%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.44%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.13%  "

# Row 44
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0673"
$ws.Range("E44").Value = "  -2.55%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.336"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.19%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0407"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.09%  "

# Row 47
$ws.Range("E47").Value = "  +4.09%  "

# Row 48
$ws.Range("E48").Value = "  -0.30%  "

# Row 49
$ws.Range("E49").Value = "  -0.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.81%  "

# Row 51
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.78"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +18.82%  "
